# "Version final de Anteproyecto"
#
# Updates the Gantt data table on Hoja1: durations/start dates for several
# activities change, which (via the existing shared formula in column D)
# also shifts the computed end dates. The stacked-bar Gantt chart reads its
# series straight from these cells, so its cached points follow along; the
# second series' fill is also restyled from the accent2 theme color to a
# neutral gray (White, Background 1, Darker 25%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "Aplicación de JavaCC en prácticas de PIAT": duration 50 -> 90 days
# (D3 = B3+C3 recalculates automatically via the shared formula)
$ws.Range("C3").Value = 90

# Row 4 - "Generación de Documentación": duration 30 -> 50 days
$ws.Range("C4").Value = 50

# Row 6 - "Preparación de la presentación del TFG": start date + duration
$ws.Range("B6").Value = 45291
$ws.Range("C6").Value = 25

# Row 7 - "Evaluación y validación del proyecto": start date + duration
$ws.Range("B7").Value = 45316
$ws.Range("C7").Value = 20

# Restyle the "Duración" series of the Gantt chart: accent2 -> bg1 (white),
# darker 25% (lumMod 75%), i.e. RGB(191,191,191).
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$durationSeries = $chart.SeriesCollection(2)
$durationSeries.Interior.Color = 12566463

# Final selected cell left on the sheet
$ws.Range("H17").Select() | Out-Null
